# Apply the price/volume (and a few coin identity) refreshes captured
# by the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.954.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.227.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.05%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.78"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.551"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0972"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "35.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.04%  "
$ws.Range("E13").Value = "  -2.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.563.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.851"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.227.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.869.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0957"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.13%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.125"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0711"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +24.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0280"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.79%  "
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.08%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.79%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.100"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.191"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.72%  "
